$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C
$ws.Range("C1").Value = "Value"

# Update the correlation values (sourced from Yahoo Finance)
$ws.Range("C2").Value = 0.674
$ws.Range("C3").Value = 0.5905
$ws.Range("C4").Value = 0.4643
$ws.Range("C5").Value = 0.3447
$ws.Range("C6").Value = -0.007
